$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "National Weather Service" record (previously row 4) moves up to row 2,
# pushing the "Man killed..." (was row 2) and "6 snowmobilers..." (was row 3)
# rows down by one position each (a cyclic rotation of rows 2-4).
#
# First clear the text/number content of rows 2-4 so the shared-string table
# drops the now-unreferenced entries; hyperlinks (attached to E2:E4) are left
# untouched by ClearContents.
$ws.Range("A2:E4").ClearContents()

# Re-populate row by row, column by column (A, B, D, E - skipping the numeric
# C column) in the exact final order so new shared strings are (re)created in
# that scan order.
$ws.Range("A2").Value = "National Weather Service"
$ws.Range("B2").Value = "2008-01-04T00:00:00UTC"
$ws.Range("D2").Value = "day_2_to_30"
$ws.Range("E2").Value = "http://www.wrh.noaa.gov/pdt/reference/20080104/index.php?wfo=pdt"

$ws.Range("A3").Value = "Man killed, seven injured on icy mountain highway"
$ws.Range("B3").Value = "2021-01-07T09:16:00UTC"
$ws.Range("D3").Value = "day_31_beyond"
$ws.Range("E3").Value = "http://www.katu.com/news/13486032.html"

$ws.Range("A4").Value = "6 snowmobilers rescued from remote outpost"
$ws.Range("B4").Value = "2008-01-07T22:11:54UTC"
$ws.Range("D4").Value = "day_2_to_30"
$ws.Range("E4").Value = "https://www.nbcnews.com/id/22539077"

# Numeric "historical distance" column follows the same rotation.
$ws.Range("C2").Value = 3
$ws.Range("C3").Value = 4755
$ws.Range("C4").Value = 6
